$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ID match")

# Insert a new row at row 10 (pushes old rows 10-18 down to 11-19, carrying
# their formatting with them so the fill-color boundary shifts too)
$ws.Rows.Item(10).Insert()

# Column A: keep the sequential numbering 1..19
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(12, 1).Value = 12
$ws.Cells.Item(13, 1).Value = 13
$ws.Cells.Item(14, 1).Value = 14
$ws.Cells.Item(15, 1).Value = 15
$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(17, 1).Value = 17
$ws.Cells.Item(18, 1).Value = 18
$ws.Cells.Item(19, 1).Value = 19

# Column B: new lookup value 9 at row 10, old 10-18 values shift down one row
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(11, 2).Value = 24
$ws.Cells.Item(12, 2).Value = 23
$ws.Cells.Item(13, 2).Value = 22
$ws.Cells.Item(14, 2).Value = 21
$ws.Cells.Item(15, 2).Value = 20
$ws.Cells.Item(16, 2).Value = 19
$ws.Cells.Item(17, 2).Value = 18
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(19, 2).Value = 16

$ws.Range("F8").Select()
